# "Bilder von den Konzepten" - add two more logged entries (concept
# sketches / design-view-model work) to the Stundenerfassung sheet and
# make it the active sheet/tab again.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Stundenerfassung")

# Copy the date formatting of the last existing row down into the two
# new rows first, so the new date cells reuse the workbook's existing
# date style (s="3") instead of creating a brand-new one.
$ws1.Range("A25").Copy()
$ws1.Range("A26:A27").PasteSpecial(-4122)  # xlPasteFormats

# Row 26: 04.06.2017 - Schriftliche Arbeit / Erstellung Skizzen - 3h
$ws1.Cells.Item(26, 1).Value = 42890
$ws1.Cells.Item(26, 2).Value = "Schriftliche Arbeit"
$ws1.Cells.Item(26, 3).Value = "Erstellung Skizzen"
$ws1.Cells.Item(26, 4).Value = 3

# Row 27: 04.06.2017 - Design View Model / Erstellung Konzept - 1h
$ws1.Cells.Item(27, 1).Value = 42890
$ws1.Cells.Item(27, 2).Value = "Design View Model"
$ws1.Cells.Item(27, 3).Value = "Erstellung Konzept"
$ws1.Cells.Item(27, 4).Value = 1

# Bring "Stundenerfassung" back to the front (it becomes the
# tabSelected / active sheet instead of "Wochen") and update the
# on-screen selection to the new bottom of the list.
$ws1.Activate()
$ws1.Range("C29").Select()
